# Add a "device_id" validation column (G) to Sheet1, mirroring the
# header/data formatting already used by the other columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from the neighboring header/data cells onto column G
# before writing values, so the new column inherits the same style
# (border + bold/centered header font) without creating new style entries.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("F2:F6").Copy()
$ws.Range("G2:G6").PasteSpecial(-4122)

# Header
$ws.Range("G1").Value = "device_id"

# Per-employee device ids. Written in this order so new shared-string
# table entries are allocated in the same sequence as the target file
# (dsf1344324 before dsf434242).
$ws.Range("G3").Value = "dsf1344324"
$ws.Range("G2").Value = "dsf434242"
$ws.Range("G4").Value = "assacds343"
$ws.Range("G5").Value = "cbfd14324"
$ws.Range("G6").Value = "vdbdf343"

# Match the width used by the other wide columns (D/E). The stored width
# is quantized by the engine to the nearest 1/6 character; 23 is the
# input that lands closest to the D/E columns' stored width.
$ws.Columns.Item(7).ColumnWidth = 23

# Move the active selection to the last filled cell in the new column
$ws.Range("G6").Select()
